# Automatic Word resume update
#
# Adds a new "Docker / Rancher container management" bullet to the
# "General skills" list, right after the existing "Continuous integration"
# bullet (same paragraph style / numbering as its neighbours).

$d = $word.ActiveDocument

# Locate the exact "Continuous integration" bullet (case-sensitive, since
# a differently-cased "Continuous Integration" bullet also exists further
# down the document, under "Specific experience").
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.Trim()
    if ($txt.Equals("Continuous integration")) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -gt 0) {
    $anchor = $d.Paragraphs.Item($targetIndex)

    # Insert a new paragraph right after the anchor; it inherits the
    # anchor's paragraph formatting (Compact style + numId 1002 bullet).
    $anchor.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($targetIndex + 1)
    $newPara.Range.InsertAfter("Docker / Rancher container management")
}
